$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSETGVANNJ%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T130003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECMaCXVzLXdlc3QtMiJGMEQCICzRdeod6FF8QxCiwwjIHfhE8X18%2FbHyT7yZZxE76qc8AiAdvQBJMZJ3gu7YPF0lRpukN1kfhH8GGXwmbvcHBANMVSrmAwiM%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM8yhRmg3Oc7gYSXj%2FKroDEI3I8WJaLvadbO1BHv%2FGztzUtTXbfTOnCu%2FuKs3d%2FeriO8UcIavpV4%2B12UBEVMrPkOdD06tiGzC%2BLa33%2F5h4G0q2xfynIM3qUn%2BFICCe4E85C7q36NRFR3WIY4NdpiFQZrXgLb%2FA0N3N1FN7cPybVeR%2FEFib0Ugd4batU9D53IdczD2m680mNQTYJHvgR64naJkFJoRQqfmPLjgM7%2F8lKtWG4hHnFn%2Fr5lgvaEOb6my0Tv5sfOQU1fRzLFpqRXR8DvsHEsWLiYhv9QHpbrPQvRRnhdcGPf2fLTg8%2Fnw7XVeabEyrtaGutajqSqbFuzeFn1DZIyrQsYZLNbEc8eBLVKi8WJiqqcGIJ9BdUhNtlJKCHmDOGjpF2YLfQL5aTuzwGLrM%2BjkacARpkphyl%2FruIjM5C3AS1j6PStEs6cPBzvHruzeYcJyrNPLjfWtR9mnAw8oXuw7qexJIfMjP73w2lyUU6agntSJnFA8hnAlVGbo3xXFMamDlwfDOlF1qKta12lVr1gHSLpnbD8IY0zXKOdwwxVvx0RFzBs7ZnVlpGUAeiQhH0ON0xI1NmpMGFbOCDCruT3IuCAUblTD1t6S%2FBjqmAcVPNUEIOpAzDtZAcWFkAGDTmhKAzq5d12gSy572kepn80PbFqLJVOBnRolvW1zgh36M1MZM40WIcCqD4hgyRNc%2Fv8Zchf8ThS3avVaLslpS4xpltlqCyWUqPAbYUQAR7aRENrPuVCQiANBvQiqUHth%2Fjy84aWwbGj1miLl%2BFpX8lW130Rz1E5JH19Lj%2B2%2FpOBOqNRFcjLDQcIjtMoH4aCl98ViBazE%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2021%3A00%3A03%20GMT&X-Amz-Signature=e021c7137b6845926aa526262430cb67fbf627a4f963885fbd694dbf83191ddd'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSETGVANNJ%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T130003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECMaCXVzLXdlc3QtMiJGMEQCICzRdeod6FF8QxCiwwjIHfhE8X18%2FbHyT7yZZxE76qc8AiAdvQBJMZJ3gu7YPF0lRpukN1kfhH8GGXwmbvcHBANMVSrmAwiM%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM8yhRmg3Oc7gYSXj%2FKroDEI3I8WJaLvadbO1BHv%2FGztzUtTXbfTOnCu%2FuKs3d%2FeriO8UcIavpV4%2B12UBEVMrPkOdD06tiGzC%2BLa33%2F5h4G0q2xfynIM3qUn%2BFICCe4E85C7q36NRFR3WIY4NdpiFQZrXgLb%2FA0N3N1FN7cPybVeR%2FEFib0Ugd4batU9D53IdczD2m680mNQTYJHvgR64naJkFJoRQqfmPLjgM7%2F8lKtWG4hHnFn%2Fr5lgvaEOb6my0Tv5sfOQU1fRzLFpqRXR8DvsHEsWLiYhv9QHpbrPQvRRnhdcGPf2fLTg8%2Fnw7XVeabEyrtaGutajqSqbFuzeFn1DZIyrQsYZLNbEc8eBLVKi8WJiqqcGIJ9BdUhNtlJKCHmDOGjpF2YLfQL5aTuzwGLrM%2BjkacARpkphyl%2FruIjM5C3AS1j6PStEs6cPBzvHruzeYcJyrNPLjfWtR9mnAw8oXuw7qexJIfMjP73w2lyUU6agntSJnFA8hnAlVGbo3xXFMamDlwfDOlF1qKta12lVr1gHSLpnbD8IY0zXKOdwwxVvx0RFzBs7ZnVlpGUAeiQhH0ON0xI1NmpMGFbOCDCruT3IuCAUblTD1t6S%2FBjqmAcVPNUEIOpAzDtZAcWFkAGDTmhKAzq5d12gSy572kepn80PbFqLJVOBnRolvW1zgh36M1MZM40WIcCqD4hgyRNc%2Fv8Zchf8ThS3avVaLslpS4xpltlqCyWUqPAbYUQAR7aRENrPuVCQiANBvQiqUHth%2Fjy84aWwbGj1miLl%2BFpX8lW130Rz1E5JH19Lj%2B2%2FpOBOqNRFcjLDQcIjtMoH4aCl98ViBazE%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2021%3A00%3A03%20GMT&X-Amz-Signature=79e300c6e50c2fa398f8f76a4a797e35547b19aa9a4e74d3d11cb8fdca9e25bf'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSETGVANNJ%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T130003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECMaCXVzLXdlc3QtMiJGMEQCICzRdeod6FF8QxCiwwjIHfhE8X18%2FbHyT7yZZxE76qc8AiAdvQBJMZJ3gu7YPF0lRpukN1kfhH8GGXwmbvcHBANMVSrmAwiM%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM8yhRmg3Oc7gYSXj%2FKroDEI3I8WJaLvadbO1BHv%2FGztzUtTXbfTOnCu%2FuKs3d%2FeriO8UcIavpV4%2B12UBEVMrPkOdD06tiGzC%2BLa33%2F5h4G0q2xfynIM3qUn%2BFICCe4E85C7q36NRFR3WIY4NdpiFQZrXgLb%2FA0N3N1FN7cPybVeR%2FEFib0Ugd4batU9D53IdczD2m680mNQTYJHvgR64naJkFJoRQqfmPLjgM7%2F8lKtWG4hHnFn%2Fr5lgvaEOb6my0Tv5sfOQU1fRzLFpqRXR8DvsHEsWLiYhv9QHpbrPQvRRnhdcGPf2fLTg8%2Fnw7XVeabEyrtaGutajqSqbFuzeFn1DZIyrQsYZLNbEc8eBLVKi8WJiqqcGIJ9BdUhNtlJKCHmDOGjpF2YLfQL5aTuzwGLrM%2BjkacARpkphyl%2FruIjM5C3AS1j6PStEs6cPBzvHruzeYcJyrNPLjfWtR9mnAw8oXuw7qexJIfMjP73w2lyUU6agntSJnFA8hnAlVGbo3xXFMamDlwfDOlF1qKta12lVr1gHSLpnbD8IY0zXKOdwwxVvx0RFzBs7ZnVlpGUAeiQhH0ON0xI1NmpMGFbOCDCruT3IuCAUblTD1t6S%2FBjqmAcVPNUEIOpAzDtZAcWFkAGDTmhKAzq5d12gSy572kepn80PbFqLJVOBnRolvW1zgh36M1MZM40WIcCqD4hgyRNc%2Fv8Zchf8ThS3avVaLslpS4xpltlqCyWUqPAbYUQAR7aRENrPuVCQiANBvQiqUHth%2Fjy84aWwbGj1miLl%2BFpX8lW130Rz1E5JH19Lj%2B2%2FpOBOqNRFcjLDQcIjtMoH4aCl98ViBazE%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2021%3A00%3A03%20GMT&X-Amz-Signature=3722ce55e011a6ef0d9ff276f6325d4505d0d56b2ea01759d868d716cefb1fb0'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSETGVANNJ%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T130003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECMaCXVzLXdlc3QtMiJGMEQCICzRdeod6FF8QxCiwwjIHfhE8X18%2FbHyT7yZZxE76qc8AiAdvQBJMZJ3gu7YPF0lRpukN1kfhH8GGXwmbvcHBANMVSrmAwiM%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM8yhRmg3Oc7gYSXj%2FKroDEI3I8WJaLvadbO1BHv%2FGztzUtTXbfTOnCu%2FuKs3d%2FeriO8UcIavpV4%2B12UBEVMrPkOdD06tiGzC%2BLa33%2F5h4G0q2xfynIM3qUn%2BFICCe4E85C7q36NRFR3WIY4NdpiFQZrXgLb%2FA0N3N1FN7cPybVeR%2FEFib0Ugd4batU9D53IdczD2m680mNQTYJHvgR64naJkFJoRQqfmPLjgM7%2F8lKtWG4hHnFn%2Fr5lgvaEOb6my0Tv5sfOQU1fRzLFpqRXR8DvsHEsWLiYhv9QHpbrPQvRRnhdcGPf2fLTg8%2Fnw7XVeabEyrtaGutajqSqbFuzeFn1DZIyrQsYZLNbEc8eBLVKi8WJiqqcGIJ9BdUhNtlJKCHmDOGjpF2YLfQL5aTuzwGLrM%2BjkacARpkphyl%2FruIjM5C3AS1j6PStEs6cPBzvHruzeYcJyrNPLjfWtR9mnAw8oXuw7qexJIfMjP73w2lyUU6agntSJnFA8hnAlVGbo3xXFMamDlwfDOlF1qKta12lVr1gHSLpnbD8IY0zXKOdwwxVvx0RFzBs7ZnVlpGUAeiQhH0ON0xI1NmpMGFbOCDCruT3IuCAUblTD1t6S%2FBjqmAcVPNUEIOpAzDtZAcWFkAGDTmhKAzq5d12gSy572kepn80PbFqLJVOBnRolvW1zgh36M1MZM40WIcCqD4hgyRNc%2Fv8Zchf8ThS3avVaLslpS4xpltlqCyWUqPAbYUQAR7aRENrPuVCQiANBvQiqUHth%2Fjy84aWwbGj1miLl%2BFpX8lW130Rz1E5JH19Lj%2B2%2FpOBOqNRFcjLDQcIjtMoH4aCl98ViBazE%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2021%3A00%3A03%20GMT&X-Amz-Signature=84df6148d75a1464b0b0163ce685fe430d704c4a8d3d2f20668e243a4ff12b08'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSETGVANNJ%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T130003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECMaCXVzLXdlc3QtMiJGMEQCICzRdeod6FF8QxCiwwjIHfhE8X18%2FbHyT7yZZxE76qc8AiAdvQBJMZJ3gu7YPF0lRpukN1kfhH8GGXwmbvcHBANMVSrmAwiM%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM8yhRmg3Oc7gYSXj%2FKroDEI3I8WJaLvadbO1BHv%2FGztzUtTXbfTOnCu%2FuKs3d%2FeriO8UcIavpV4%2B12UBEVMrPkOdD06tiGzC%2BLa33%2F5h4G0q2xfynIM3qUn%2BFICCe4E85C7q36NRFR3WIY4NdpiFQZrXgLb%2FA0N3N1FN7cPybVeR%2FEFib0Ugd4batU9D53IdczD2m680mNQTYJHvgR64naJkFJoRQqfmPLjgM7%2F8lKtWG4hHnFn%2Fr5lgvaEOb6my0Tv5sfOQU1fRzLFpqRXR8DvsHEsWLiYhv9QHpbrPQvRRnhdcGPf2fLTg8%2Fnw7XVeabEyrtaGutajqSqbFuzeFn1DZIyrQsYZLNbEc8eBLVKi8WJiqqcGIJ9BdUhNtlJKCHmDOGjpF2YLfQL5aTuzwGLrM%2BjkacARpkphyl%2FruIjM5C3AS1j6PStEs6cPBzvHruzeYcJyrNPLjfWtR9mnAw8oXuw7qexJIfMjP73w2lyUU6agntSJnFA8hnAlVGbo3xXFMamDlwfDOlF1qKta12lVr1gHSLpnbD8IY0zXKOdwwxVvx0RFzBs7ZnVlpGUAeiQhH0ON0xI1NmpMGFbOCDCruT3IuCAUblTD1t6S%2FBjqmAcVPNUEIOpAzDtZAcWFkAGDTmhKAzq5d12gSy572kepn80PbFqLJVOBnRolvW1zgh36M1MZM40WIcCqD4hgyRNc%2Fv8Zchf8ThS3avVaLslpS4xpltlqCyWUqPAbYUQAR7aRENrPuVCQiANBvQiqUHth%2Fjy84aWwbGj1miLl%2BFpX8lW130Rz1E5JH19Lj%2B2%2FpOBOqNRFcjLDQcIjtMoH4aCl98ViBazE%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2021%3A00%3A03%20GMT&X-Amz-Signature=7e09135dd8de8aaac9ad7329cdd08041e7f4ca7d1116146a2d499595a4d13e5b'
